$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $addr, $text)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-CellText $ws 'D2' '42.780.15'
Set-CellText $ws 'E2' '  -1.30%  '
Set-CellText $ws 'D3' '2.354.01'
Set-CellText $ws 'E3' '  -1.72%  '
Set-CellText $ws 'E4' '  -0.11%  '
Set-CellText $ws 'D5' '318.97'
Set-CellText $ws 'E5' '  -1.24%  '
Set-CellText $ws 'D6' '106.78'
Set-CellText $ws 'E6' '  +1.31%  '
Set-CellText $ws 'E7' '  -2.58%  '
Set-CellText $ws 'E8' '  +0.01%  '
Set-CellText $ws 'D9' '0.621'
Set-CellText $ws 'E9' '  -5.03%  '
Set-CellText $ws 'D10' '41.49'
Set-CellText $ws 'E10' '  -1.07%  '
Set-CellText $ws 'D11' '0.0926'
Set-CellText $ws 'E11' '  -1.97%  '
Set-CellText $ws 'D12' '8.45'
Set-CellText $ws 'E12' '  -1.91%  '
Set-CellText $ws 'E13' '  -1.57%  '
Set-CellText $ws 'E14' '  -0.11%  '
Set-CellText $ws 'D15' '16.00'
Set-CellText $ws 'E15' '  -9.03%  '
Set-CellText $ws 'D16' '2.708.85'
Set-CellText $ws 'E16' '  -1.62%  '
Set-CellText $ws 'D17' '2.336.66'
Set-CellText $ws 'E17' '  -2.41%  '
Set-CellText $ws 'D18' '42.731.67'
Set-CellText $ws 'E18' '  -1.43%  '
Set-CellText $ws 'E19' '  +4.31%  '
Set-CellText $ws 'E20' '  -2.47%  '
Set-CellText $ws 'D21' '77.13'
Set-CellText $ws 'E21' '  +0.88%  '
Set-CellText $ws 'E22' '  +5.13%  '
Set-CellText $ws 'D23' '258.98'
Set-CellText $ws 'E23' '  -4.65%  '
Set-CellText $ws 'E24' '  -5.05%  '
Set-CellText $ws 'D25' '9.44'
Set-CellText $ws 'E25' '  -5.15%  '
Set-CellText $ws 'E26' '  +0.01%  '
Set-CellText $ws 'D27' '11.41'
Set-CellText $ws 'E27' '  -3.72%  '
Set-CellText $ws 'E28' '  +0.15%  '
Set-CellText $ws 'E29' '  +0.89%  '
Set-CellText $ws 'D30' '175.12'
Set-CellText $ws 'E30' '  -1.12%  '
Set-CellText $ws 'E31' '  -3.43%  '
Set-CellText $ws 'B32' 'Filecoin'
Set-CellText $ws 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText $ws 'D32' '6.14'
Set-CellText $ws 'E32' '  +4.02%  '
Set-CellText $ws 'B33' 'Hedera'
Set-CellText $ws 'C33' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText $ws 'D33' '0.0886'
Set-CellText $ws 'E33' '  -4.68%  '
Set-CellText $ws 'E34' '  -8.01%  '
Set-CellText $ws 'E35' '  +13.81%  '
Set-CellText $ws 'E36' '  -2.73%  '
Set-CellText $ws 'E37' '  -5.35%  '
Set-CellText $ws 'E38' '  -1.94%  '
Set-CellText $ws 'E39' '  -7.43%  '
Set-CellText $ws 'E40' '  -4.89%  '
Set-CellText $ws 'D41' '0.239'
Set-CellText $ws 'E41' '  +1.89%  '
Set-CellText $ws 'D42' '71.85'
Set-CellText $ws 'E42' '  +3.37%  '
Set-CellText $ws 'E43' '  -8.22%  '
Set-CellText $ws 'E44' '  -0.18%  '
Set-CellText $ws 'D45' '114.24'
Set-CellText $ws 'E45' '  -10.29%  '
Set-CellText $ws 'D46' '12.00'
Set-CellText $ws 'E46' '  -4.91%  '
Set-CellText $ws 'D47' '5.49'
Set-CellText $ws 'E47' '  -2.93%  '
Set-CellText $ws 'D48' '9.13'
Set-CellText $ws 'E48' '  -5.18%  '
Set-CellText $ws 'D49' '85.06'
Set-CellText $ws 'E49' '  -4.70%  '
Set-CellText $ws 'D50' '73.93'
Set-CellText $ws 'E50' '  +1.28%  '
Set-CellText $ws 'E51' '  -1.49%  '
